$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 568, shifting the existing rows 568-609 down to 569-610.
$ws.Rows.Item(568).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a date-like string ("2026/01/06"); format it as Text first so it is
# stored as a literal string rather than being auto-converted to a date serial number,
# then clear the formatting so the cell ends up with the default (unstyled) look,
# matching the rest of the column.
$ws.Range("A568").NumberFormat = "@"
$ws.Range("A568").Value = "2026/01/06"
$ws.Range("A568").ClearFormats()

$ws.Range("B568").Value = "火"
$ws.Range("C568").Value = 19
$ws.Range("D568").Value = 159
